# Updated TPM-derived NATMI ligand-receptor values (rows 2-17, columns G-J, M-T).
# New expression levels change ligand/receptor average & total expression,
# their derived-specificity fractions, and the resulting edge weights/specificities.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.082188
$ws.Range("H2").Value = 24.246564
$ws.Range("I2").Value = 0.2755301789948819
$ws.Range("J2").Value = 0.2755301789948819
$ws.Range("M2").Value = 57.48524866666667
$ws.Range("N2").Value = 172.455746
$ws.Range("O2").Value = 0.2249897472933344
$ws.Range("P2").Value = 0.2249897472933344
$ws.Range("Q2").Value = 464.6065869507494
$ws.Range("R2").Value = 4181.459282556744
$ws.Range("S2").Value = 0.06199146534374567
$ws.Range("T2").Value = 0.06199146534374567

# Row 3
$ws.Range("G3").Value = 8.082188
$ws.Range("H3").Value = 24.246564
$ws.Range("I3").Value = 0.2755301789948819
$ws.Range("J3").Value = 0.2755301789948819
$ws.Range("M3").Value = 72.97955566666666
$ws.Range("N3").Value = 218.938667
$ws.Range("O3").Value = 0.2856324390668287
$ws.Range("P3").Value = 0.2856324390668287
$ws.Range("Q3").Value = 589.8344890544653
$ws.Range("R3").Value = 5308.510401490188
$ws.Range("S3").Value = 0.07870035706282799
$ws.Range("T3").Value = 0.07870035706282801

# Row 4
$ws.Range("G4").Value = 8.082188
$ws.Range("H4").Value = 24.246564
$ws.Range("I4").Value = 0.2755301789948819
$ws.Range("J4").Value = 0.2755301789948819
$ws.Range("M4").Value = 50.94830300000001
$ws.Range("N4").Value = 152.844909
$ws.Range("O4").Value = 0.1994049966359642
$ws.Range("P4").Value = 0.1994049966359642
$ws.Range("Q4").Value = 411.7737631269641
$ws.Range("R4").Value = 3705.963868142677
$ws.Range("S4").Value = 0.05494209441558103
$ws.Range("T4").Value = 0.05494209441558104

# Row 5
$ws.Range("G5").Value = 8.082188
$ws.Range("H5").Value = 24.246564
$ws.Range("I5").Value = 0.2755301789948819
$ws.Range("J5").Value = 0.2755301789948819
$ws.Range("M5").Value = 74.08852933333333
$ws.Range("N5").Value = 222.265588
$ws.Range("O5").Value = 0.2899728170038728
$ws.Range("P5").Value = 0.2899728170038728
$ws.Range("Q5").Value = 598.7974227155146
$ws.Range("R5").Value = 5389.176804439632
$ws.Range("S5").Value = 0.07989626217272719
$ws.Range("T5").Value = 0.0798962621727272

# Row 6
$ws.Range("G6").Value = 13.29805733333333
$ws.Range("H6").Value = 39.894172
$ws.Range("I6").Value = 0.4533445791334642
$ws.Range("J6").Value = 0.4533445791334642
$ws.Range("M6").Value = 57.48524866666667
$ws.Range("N6").Value = 172.455746
$ws.Range("O6").Value = 0.2249897472933344
$ws.Range("P6").Value = 0.2249897472933344
$ws.Range("Q6").Value = 764.4421325902568
$ws.Range("R6").Value = 6879.979193312312
$ws.Range("S6").Value = 0.1019978822960411
$ws.Range("T6").Value = 0.1019978822960411

# Row 7
$ws.Range("G7").Value = 13.29805733333333
$ws.Range("H7").Value = 39.894172
$ws.Range("I7").Value = 0.4533445791334642
$ws.Range("J7").Value = 0.4533445791334642
$ws.Range("M7").Value = 72.97955566666666
$ws.Range("N7").Value = 218.938667
$ws.Range("O7").Value = 0.2856324390668287
$ws.Range("P7").Value = 0.2856324390668287
$ws.Range("Q7").Value = 970.4863154165247
$ws.Range("R7").Value = 8734.376838748723
$ws.Range("S7").Value = 0.1294899178756163
$ws.Range("T7").Value = 0.1294899178756163

# Row 8
$ws.Range("G8").Value = 13.29805733333333
$ws.Range("H8").Value = 39.894172
$ws.Range("I8").Value = 0.4533445791334642
$ws.Range("J8").Value = 0.4533445791334642
$ws.Range("M8").Value = 50.94830300000001
$ws.Range("N8").Value = 152.844909
$ws.Range("O8").Value = 0.1994049966359642
$ws.Range("P8").Value = 0.1994049966359642
$ws.Range("Q8").Value = 677.5134543300387
$ws.Range("R8").Value = 6097.621088970349
$ws.Range("S8").Value = 0.09039917427704103
$ws.Range("T8").Value = 0.09039917427704103

# Row 9
$ws.Range("G9").Value = 13.29805733333333
$ws.Range("H9").Value = 39.894172
$ws.Range("I9").Value = 0.4533445791334642
$ws.Range("J9").Value = 0.4533445791334642
$ws.Range("M9").Value = 74.08852933333333
$ws.Range("N9").Value = 222.265588
$ws.Range("O9").Value = 0.2899728170038728
$ws.Range("P9").Value = 0.2899728170038728
$ws.Range("Q9").Value = 985.2335108170149
$ws.Range("R9").Value = 8867.101597353134
$ws.Range("S9").Value = 0.1314576046847657
$ws.Range("T9").Value = 0.1314576046847657

# Row 10
$ws.Range("G10").Value = 5.789497666666667
$ws.Range("H10").Value = 17.368493
$ws.Range("I10").Value = 0.1973699855023315
$ws.Range("J10").Value = 0.1973699855023315
$ws.Range("M10").Value = 57.48524866666667
$ws.Range("N10").Value = 172.455746
$ws.Range("O10").Value = 0.2249897472933344
$ws.Range("P10").Value = 0.2249897472933344
$ws.Range("Q10").Value = 332.8107130234198
$ws.Range("R10").Value = 2995.296417210778
$ws.Range("S10").Value = 0.04440622316145863
$ws.Range("T10").Value = 0.04440622316145863

# Row 11
$ws.Range("G11").Value = 5.789497666666667
$ws.Range("H11").Value = 17.368493
$ws.Range("I11").Value = 0.1973699855023315
$ws.Range("J11").Value = 0.1973699855023315
$ws.Range("M11").Value = 72.97955566666666
$ws.Range("N11").Value = 218.938667
$ws.Range("O11").Value = 0.2856324390668287
$ws.Range("P11").Value = 0.2856324390668287
$ws.Range("Q11").Value = 422.5149672465367
$ws.Range("R11").Value = 3802.634705218831
$ws.Range("S11").Value = 0.05637527035761555
$ws.Range("T11").Value = 0.05637527035761556

# Row 12
$ws.Range("G12").Value = 5.789497666666667
$ws.Range("H12").Value = 17.368493
$ws.Range("I12").Value = 0.1973699855023315
$ws.Range("J12").Value = 0.1973699855023315
$ws.Range("M12").Value = 50.94830300000001
$ws.Range("N12").Value = 152.844909
$ws.Range("O12").Value = 0.1994049966359642
$ws.Range("P12").Value = 0.1994049966359642
$ws.Range("Q12").Value = 294.9650813391264
$ws.Range("R12").Value = 2654.685732052138
$ws.Range("S12").Value = 0.03935656129513271
$ws.Range("T12").Value = 0.03935656129513272

# Row 13
$ws.Range("G13").Value = 5.789497666666667
$ws.Range("H13").Value = 17.368493
$ws.Range("I13").Value = 0.1973699855023315
$ws.Range("J13").Value = 0.1973699855023315
$ws.Range("M13").Value = 74.08852933333333
$ws.Range("N13").Value = 222.265588
$ws.Range("O13").Value = 0.2899728170038728
$ws.Range("P13").Value = 0.2899728170038728
$ws.Range("Q13").Value = 428.9353677020982
$ws.Range("R13").Value = 3860.418309318884
$ws.Range("S13").Value = 0.05723193068812459
$ws.Range("T13").Value = 0.0572319306881246

# Row 14
$ws.Range("G14").Value = 2.163479333333334
$ws.Range("H14").Value = 6.490438
$ws.Range("I14").Value = 0.0737552563693224
$ws.Range("J14").Value = 0.0737552563693224
$ws.Range("M14").Value = 57.48524866666667
$ws.Range("N14").Value = 172.455746
$ws.Range("O14").Value = 0.2249897472933344
$ws.Range("P14").Value = 0.2249897472933344
$ws.Range("Q14").Value = 124.3681474618609
$ws.Range("R14").Value = 1119.313327156748
$ws.Range("S14").Value = 0.01659417649208894
$ws.Range("T14").Value = 0.01659417649208894

# Row 15
$ws.Range("G15").Value = 2.163479333333334
$ws.Range("H15").Value = 6.490438
$ws.Range("I15").Value = 0.0737552563693224
$ws.Range("J15").Value = 0.0737552563693224
$ws.Range("M15").Value = 72.97955566666666
$ws.Range("N15").Value = 218.938667
$ws.Range("O15").Value = 0.2856324390668287
$ws.Range("P15").Value = 0.2856324390668287
$ws.Range("Q15").Value = 157.8897604406829
$ws.Range("R15").Value = 1421.007843966146
$ws.Range("S15").Value = 0.02106689377076881
$ws.Range("T15").Value = 0.02106689377076881

# Row 16
$ws.Range("G16").Value = 2.163479333333334
$ws.Range("H16").Value = 6.490438
$ws.Range("I16").Value = 0.0737552563693224
$ws.Range("J16").Value = 0.0737552563693224
$ws.Range("M16").Value = 50.94830300000001
$ws.Range("N16").Value = 152.844909
$ws.Range("O16").Value = 0.1994049966359642
$ws.Range("P16").Value = 0.1994049966359642
$ws.Range("Q16").Value = 110.2256006089047
$ws.Range("R16").Value = 992.0304054801422
$ws.Range("S16").Value = 0.01470716664820941
$ws.Range("T16").Value = 0.01470716664820941

# Row 17
$ws.Range("G17").Value = 2.163479333333334
$ws.Range("H17").Value = 6.490438
$ws.Range("I17").Value = 0.0737552563693224
$ws.Range("J17").Value = 0.0737552563693224
$ws.Range("M17").Value = 74.08852933333333
$ws.Range("N17").Value = 222.265588
$ws.Range("O17").Value = 0.2899728170038728
$ws.Range("P17").Value = 0.2899728170038728
$ws.Range("Q17").Value = 160.2890020497271
$ws.Range("R17").Value = 1442.601018447544
$ws.Range("S17").Value = 0.02138701945825525
$ws.Range("T17").Value = 0.02138701945825525
